# "Generate Report for Handback" — fill in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns (F, G, H) for both
# language sheets, and flip Status from "Ready for handoff" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor = 15570276       # BGR encoding of RGB 6495ED (the workbook's HyperLink font colour)

function Set-HandbackRow {
    param($ws, $row, $sourceMdName, $sourceMdUrl, $targetXlfName, $targetXlfUrl, $handbackDateTime)

    # Status -> handed back
    $ws.Cells.Item($row, 3).Value = "Handed back: in sync with en-US"

    # F = Latest Target File (same file that was handed off, round-tripped back)
    $fCell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($fCell, $sourceMdUrl, "", "", $sourceMdName)
    $fCell.Font.Underline = $hyperlinkUnderline
    $fCell.Font.Color = $hyperlinkColor

    # G = Latest Handback File (the translated xlf that came back)
    $gCell = $ws.Cells.Item($row, 7)
    $ws.Hyperlinks.Add($gCell, $targetXlfUrl, "", "", $targetXlfName)
    $gCell.Font.Underline = $hyperlinkUnderline
    $gCell.Font.Color = $hyperlinkColor

    # H = Latest Handback DateTime
    $ws.Cells.Item($row, 8).Value = $handbackDateTime
}

# ---------------- Overview sheet ----------------
# The Overview sheet's zh-cn/de-de status columns mirror the same
# "Status" text as the per-language sheets, so they flip too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Set-HandbackRow `
    $wsZh `
    2 `
    "8fc7b040-dee2-421b-8e4f-1d316658501e.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a16392e0397e499a2de2d0049aa07717f8fc03d7/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md" `
    "8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57f31e446e4a9bd08c216550d6136faaaa3ba59f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.zh-cn.xlf" `
    "2016-03-19 12:35:23"

Set-HandbackRow `
    $wsZh `
    3 `
    "ba5b2a27-c42e-4541-af27-6057b6515bf4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a16392e0397e499a2de2d0049aa07717f8fc03d7/e2e/ba5b2a27-c42e-4541-af27-6057b6515bf4.md" `
    "ba5b2a27-c42e-4541-af27-6057b6515bf4.ea2d99e310fa7743fbd89285d41bf9b9d1b88fed.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/57f31e446e4a9bd08c216550d6136faaaa3ba59f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ba5b2a27-c42e-4541-af27-6057b6515bf4.ea2d99e310fa7743fbd89285d41bf9b9d1b88fed.zh-cn.xlf" `
    "2016-03-19 12:35:23"

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")

Set-HandbackRow `
    $wsDe `
    2 `
    "8fc7b040-dee2-421b-8e4f-1d316658501e.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a16392e0397e499a2de2d0049aa07717f8fc03d7/e2e/8fc7b040-dee2-421b-8e4f-1d316658501e.md" `
    "8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/54d2fe8ad8de10b40ef4a54b6f94b3e29d54d45f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8fc7b040-dee2-421b-8e4f-1d316658501e.202cae47a453ea5feb7fbbec71dc0f8dbb40a093.de-de.xlf" `
    "2016-03-19 12:35:29"

Set-HandbackRow `
    $wsDe `
    3 `
    "ba5b2a27-c42e-4541-af27-6057b6515bf4.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a16392e0397e499a2de2d0049aa07717f8fc03d7/e2e/ba5b2a27-c42e-4541-af27-6057b6515bf4.md" `
    "ba5b2a27-c42e-4541-af27-6057b6515bf4.ea2d99e310fa7743fbd89285d41bf9b9d1b88fed.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/54d2fe8ad8de10b40ef4a54b6f94b3e29d54d45f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ba5b2a27-c42e-4541-af27-6057b6515bf4.ea2d99e310fa7743fbd89285d41bf9b9d1b88fed.de-de.xlf" `
    "2016-03-19 12:35:29"
